$d = $word.ActiveDocument

function SetParaText($index, $new) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    $r.End = $r.End - 1
    $r.Text = $new
}

SetParaText 1 "ContosoLearn Market Research"

SetParaText 2 "AdatumLearn: AdatumLearn is a top AI-powered learning platform that uses artificial intelligence to enrich eLearning with features that automate a variety of tasks. It is known for its content authoring capabilities and adaptive learning technology."

SetParaText 3 "AdventureLearn: AdventureLearn is another AI-powered learning platform that offers personalized learning experiences and data-driven recommendations."

SetParaText 4 "AlpineTraining: AlpineTraining is a mobile-first learning platform that focuses on microlearning."

SetParaText 5 "Bellows OnDemand: Bellows OnDemand is a comprehensive learning solution that offers content creation and social collaboration."

SetParaText 6 "FabrikamLearning: FabrikamLearning provides a suite of learning platforms that cater to different learning needs."

SetParaText 7 "FirstUp Cards: FirstUp Cards is a mobile learning app that is ideal for training on safety procedures, compliance, new product knowledge or any other type of training scenario."

SetParaText 8 "Munson'sLearn: Munson'sLearn is designed to enable businesses to train their employees, partners, and customers."

SetParaText 9 "LibertyLearn: LibertyLearn is a fast LMS for your mission-critical project."

SetParaText 10 "WoodgroveLMS: WoodgroveLMS is a functional and attractive learning management system built to provide a best-in-class training experience."

SetParaText 11 "NorthwindWorlds: NorthwindWorlds is a powerful, easy-to-use, and reliable training solution for individuals and enterprises."

SetParaText 12 "ProsewareLearn: ProsewareLearn is an online education company that offers a variety of video training courses for software developers, IT administrators, and creative professionals through its website."

SetParaText 13 "RelecloudLearn: RelecloudLearn is an American online learning platform that offers massive open online courses (MOOC), specializations, and degrees in a variety of subjects."

SetParaText 14 "TreyAcademy: TreyAcademy is an online learning platform aimed at professional adults and students, developed in May 2010."

SetParaText 15 "These platforms have a significant market presence and are widely recognized for their AI-powered features, such as personalized learning experiences, data-driven recommendations, and automation of tasks. They are transforming the eLearning landscape by leveraging AI to deliver more engaging, rewarding, and personalized learning experiences. "

Write-Output "Done"
